# Edit script: updates the 9 remaining vocabulary slides with new kanji/
# definitions/page ranges, and removes the trailing 5 slides (10-14) that
# were dropped from the deck.

$p = $ppt.ActivePresentation

# --- Update slide 1: 境 -> 殺す ---------------------------------------
$s = $p.Slides.Item(1)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "殺す"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "to kill, to slay, to murder, to slaughter | to suppress, to block, to hamper, to destroy (e.g. talent), to eliminate (e.g..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "67-68"

# --- Update slide 2: 環境 -> 殺人 --------------------------------------
$s = $p.Slides.Item(2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "殺人"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "murder, homicide, manslaughter..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "67-68"

# --- Update slide 3: 音響 -> 農薬 --------------------------------------
$s = $p.Slides.Item(3)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "農薬"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "agricultural chemical (i.e. pesticide, herbicide, fungicide, etc.), agrochemical, agrichemical..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "67-68"

# --- Update slide 4: 影 -> 収入印紙 ------------------------------------
$s = $p.Slides.Item(4)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "収入印紙"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "revenue stamp..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "67-68"

# --- Update slide 5: 影響 -> 収穫 --------------------------------------
$s = $p.Slides.Item(5)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "収穫"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "harvest, crop, ingathering | fruits (of one's labors), gain, result, returns..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "67-68"

# --- Update slide 6: 感染 -> 少量 --------------------------------------
$s = $p.Slides.Item(6)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "少量"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "small quantity, small amount | narrowmindedness..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "67-68"

# --- Update slide 7: 汚い -> 完全 --------------------------------------
$s = $p.Slides.Item(7)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "完全"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "perfect, complete..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "67-68"

# --- Update slide 8: 汚す -> 原因 --------------------------------------
$s = $p.Slides.Item(8)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "原因"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "cause, origin, source..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "67-68"

# --- Update slide 9: 汚染 -> 一環 --------------------------------------
$s = $p.Slides.Item(9)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "一環"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "link (e.g. in a chain of events), part (of a plan, campaign, activities, etc.) | monocyclic..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "67-68"

# --- Remove the trailing 5 slides (old slides 10-14) -------------------
$p.Slides.Item(10).Delete()
$p.Slides.Item(10).Delete()
$p.Slides.Item(10).Delete()
$p.Slides.Item(10).Delete()
$p.Slides.Item(10).Delete()
